$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1586826347305389
$ws.Range("C2").Value = 0.6347305389221557
$ws.Range("J2").Value = 0.02395209580838323
$ws.Range("P2").Value = 0.1167664670658683
$ws.Range("S2").Value = 0.0658682634730539
$ws.Range("B3").Value = 0.004651162790697674
$ws.Range("C3").Value = 0.009302325581395349
$ws.Range("J3").Value = 0.0186046511627907
$ws.Range("P3").Value = 0.7767441860465116
$ws.Range("S3").Value = 0.1906976744186047
$ws.Range("J4").Value = 0.07692307692307693
$ws.Range("P4").Value = 0.5384615384615384
$ws.Range("S4").Value = 0.3846153846153846
$ws.Range("S5").Value = 1
$ws.Range("B6").Value = 0.07514450867052024
$ws.Range("F6").Value = 0.03468208092485549
$ws.Range("J6").Value = 0.3121387283236994
$ws.Range("O6").Value = 0.01734104046242774
$ws.Range("Q6").Value = 0.1849710982658959
$ws.Range("R6").Value = 0.06358381502890173
$ws.Range("S6").Value = 0.3121387283236994
$ws.Range("B7").Value = 0.1104651162790698
$ws.Range("D7").Value = 0.01162790697674419
$ws.Range("F7").Value = 0.04651162790697674
$ws.Range("J7").Value = 0.1337209302325581
$ws.Range("O7").Value = 0.005813953488372093
$ws.Range("Q7").Value = 0.2383720930232558
$ws.Range("R7").Value = 0.05813953488372093
$ws.Range("S7").Value = 0.3953488372093023
$ws.Range("B8").Value = 0.1278538812785388
$ws.Range("D8").Value = 0.0091324200913242
$ws.Range("E8").Value = 0.00228310502283105
$ws.Range("F8").Value = 0.0547945205479452
$ws.Range("J8").Value = 0.1324200913242009
$ws.Range("O8").Value = 0.02511415525114155
$ws.Range("Q8").Value = 0.1757990867579909
$ws.Range("R8").Value = 0.0776255707762557
$ws.Range("S8").Value = 0.3949771689497717
$ws.Range("B9").Value = 0.125
$ws.Range("D9").Value = 0.02717391304347826
$ws.Range("F9").Value = 0.03260869565217391
$ws.Range("J9").Value = 0.1358695652173913
$ws.Range("O9").Value = 0.0108695652173913
$ws.Range("Q9").Value = 0.2010869565217391
$ws.Range("R9").Value = 0.05978260869565218
$ws.Range("S9").Value = 0.4076086956521739
$ws.Range("B10").Value = 0.1332250203086921
$ws.Range("D10").Value = 0.02274573517465475
$ws.Range("F10").Value = 0.06417546709991877
$ws.Range("J10").Value = 0.1502843216896832
$ws.Range("O10").Value = 0.01299756295694557
$ws.Range("Q10").Value = 0.2282696994313566
$ws.Range("R10").Value = 0.06742485783915515
$ws.Range("S10").Value = 0.3208773354995939
$ws.Range("G11").Value = 0.1333333333333333
$ws.Range("J11").Value = 0.09019607843137255
$ws.Range("K11").Value = 0.2
$ws.Range("L11").Value = 0.5607843137254902
$ws.Range("S11").Value = 0.01568627450980392
$ws.Range("G12").Value = 0.7483443708609272
$ws.Range("J12").Value = 0.1854304635761589
$ws.Range("L12").Value = 0.04635761589403974
$ws.Range("S12").Value = 0.01986754966887417
$ws.Range("G13").Value = 0.6904761904761905
$ws.Range("J13").Value = 0.2619047619047619
$ws.Range("S13").Value = 0.04761904761904762
$ws.Range("F15").Value = 0.01851851851851852
$ws.Range("H15").Value = 0.1481481481481481
$ws.Range("I15").Value = 0.1172839506172839
$ws.Range("J15").Value = 0.3518518518518519
$ws.Range("K15").Value = 0.03703703703703703
$ws.Range("M15").Value = 0.006172839506172839
$ws.Range("O15").Value = 0.04320987654320987
$ws.Range("S15").Value = 0.2777777777777778
$ws.Range("F16").Value = 0.01408450704225352
$ws.Range("H16").Value = 0.2018779342723005
$ws.Range("I16").Value = 0.107981220657277
$ws.Range("J16").Value = 0.4084507042253521
$ws.Range("K16").Value = 0.07511737089201878
$ws.Range("M16").Value = 0.0187793427230047
$ws.Range("O16").Value = 0.05633802816901409
$ws.Range("S16").Value = 0.1173708920187793
$ws.Range("F17").Value = 0.01511879049676026
$ws.Range("H17").Value = 0.2116630669546436
$ws.Range("I17").Value = 0.08855291576673865
$ws.Range("J17").Value = 0.4038876889848812
$ws.Range("K17").Value = 0.07775377969762419
$ws.Range("M17").Value = 0.01727861771058315
$ws.Range("O17").Value = 0.03455723542116631
$ws.Range("S17").Value = 0.1511879049676026
$ws.Range("H18").Value = 0.1891891891891892
$ws.Range("I18").Value = 0.1081081081081081
$ws.Range("J18").Value = 0.4256756756756757
$ws.Range("K18").Value = 0.1148648648648649
$ws.Range("M18").Value = 0.006756756756756757
$ws.Range("O18").Value = 0.03378378378378379
$ws.Range("S18").Value = 0.1216216216216216
$ws.Range("F19").Value = 0.009683098591549295
$ws.Range("H19").Value = 0.221830985915493
$ws.Range("I19").Value = 0.07570422535211267
$ws.Range("J19").Value = 0.3741197183098591
$ws.Range("K19").Value = 0.1117957746478873
$ws.Range("M19").Value = 0.0272887323943662
$ws.Range("N19").Value = 0.00176056338028169
